{"js": "const pairs = [\n  [\"2024-07-14 Sunday\", \"2024-07-15 Monday\"],\n  [\"430\u00d72=\", \"780\u00d78=\"],\n  [\"943\u00d72=\", \"521\u00d75=\"],\n  [\"647\u00d78=\", \"574\u00d76=\"],\n  [\"741\u00d75=\", \"403\u00d75=\"],\n  [\"352\u00d75=\", \"401\u00d79=\"],\n  [\"617\u00d73=\", \"763\u00d79=\"],\n  [\"598\u00d73=\", \"844\u00d79=\"],\n  [\"314\u00d72=\", \"480\u00d75=\"],\n  [\"504\u00d73=\", \"305\u00d79=\"],\n  [\"144\u00d72=\", \"700\u00d73=\"],\n  [\"861\u00d73=\", \"762\u00d78=\"],\n  [\"551\u00d73=\", \"415\u00d77=\"],\n  [\"696\u00d75=\", \"358\u00d74=\"],\n  [\"942\u00d75=\", \"671\u00d73=\"],\n  [\"388\u00d78=\", \"572\u00d77=\"],\n  [\"945\u00d76=\", \"229\u00d73=\"],\n  [\"159\u00d73=\", \"462\u00d79=\"],\n  [\"760\u00d76=\", \"114\u00d76=\"],\n  [\"242\u00d74=\", \"244\u00d78=\"],\n  [\"769\u00d73=\", \"769\u00d78=\"],\n  [\"191\u00d78=\", \"621\u00d79=\"],\n  [\"811\u00d78=\", \"436\u00d78=\"],\n  [\"777\u00d73=\", \"895\u00d79=\"],\n  [\"705\u00d72=\", \"982\u00d75=\"],\n  [\"534\u00d77=\", \"599\u00d76=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$found = $find.Execute(\"2024-07-14 Sunday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-07-15 Monday\", 2)\nif (-not $found) { throw \"Not found: 2024-07-14 Sunday\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"430\u00d72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"780\u00d78=\", 2)\nif (-not $found) { throw \"Not found: 430\u00d72=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"943\u00d72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"521\u00d75=\", 2)\nif (-not $found) { throw \"Not found: 943\u00d72=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"647\u00d78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"574\u00d76=\", 2)\nif (-not $found) { throw \"Not found: 647\u00d78=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"741\u00d75=\", $false, $false, $false, $false, $false, $true, 1, $false, \"403\u00d75=\", 2)\nif (-not $found) { throw \"Not found: 741\u00d75=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"352\u00d75=\", $false, $false, $false, $false, $false, $true, 1, $false, \"401\u00d79=\", 2)\nif (-not $found) { throw \"Not found: 352\u00d75=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"617\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"763\u00d79=\", 2)\nif (-not $found) { throw \"Not found: 617\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"598\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"844\u00d79=\", 2)\nif (-not $found) { throw \"Not found: 598\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"314\u00d72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"480\u00d75=\", 2)\nif (-not $found) { throw \"Not found: 314\u00d72=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"504\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"305\u00d79=\", 2)\nif (-not $found) { throw \"Not found: 504\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"144\u00d72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"700\u00d73=\", 2)\nif (-not $found) { throw \"Not found: 144\u00d72=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"861\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"762\u00d78=\", 2)\nif (-not $found) { throw \"Not found: 861\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"551\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"415\u00d77=\", 2)\nif (-not $found) { throw \"Not found: 551\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"696\u00d75=\", $false, $false, $false, $false, $false, $true, 1, $false, \"358\u00d74=\", 2)\nif (-not $found) { throw \"Not found: 696\u00d75=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"942\u00d75=\", $false, $false, $false, $false, $false, $true, 1, $false, \"671\u00d73=\", 2)\nif (-not $found) { throw \"Not found: 942\u00d75=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"388\u00d78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"572\u00d77=\", 2)\nif (-not $found) { throw \"Not found: 388\u00d78=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"945\u00d76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"229\u00d73=\", 2)\nif (-not $found) { throw \"Not found: 945\u00d76=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"159\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"462\u00d79=\", 2)\nif (-not $found) { throw \"Not found: 159\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"760\u00d76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"114\u00d76=\", 2)\nif (-not $found) { throw \"Not found: 760\u00d76=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"242\u00d74=\", $false, $false, $false, $false, $false, $true, 1, $false, \"244\u00d78=\", 2)\nif (-not $found) { throw \"Not found: 242\u00d74=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"769\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"769\u00d78=\", 2)\nif (-not $found) { throw \"Not found: 769\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"191\u00d78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"621\u00d79=\", 2)\nif (-not $found) { throw \"Not found: 191\u00d78=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"811\u00d78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"436\u00d78=\", 2)\nif (-not $found) { throw \"Not found: 811\u00d78=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"777\u00d73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"895\u00d79=\", 2)\nif (-not $found) { throw \"Not found: 777\u00d73=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"705\u00d72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"982\u00d75=\", 2)\nif (-not $found) { throw \"Not found: 705\u00d72=\" }\n\n$find = $d.Content.Find\n$found = $find.Execute(\"534\u00d77=\", $false, $false, $false, $false, $false, $true, 1, $false, \"599\u00d76=\", 2)\nif (-not $found) { throw \"Not found: 534\u00d77=\" }\n"}
